# Updates cryptocurrency price/volume data (columns D and E) for rows 2-51
# per the latest scrape, matching the GitHub Actions commit.

function Set-CellText {
    param($ws, $addr, $text)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" '67.877.69'
Set-CellText $ws "E2" '  +0.70%  '
Set-CellText $ws "D3" '3.520.71'
Set-CellText $ws "E3" '  +0.71%  '
Set-CellText $ws "E4" '  +0.20%  '
Set-CellText $ws "D5" '610.40'
Set-CellText $ws "E5" '  +1.26%  '
Set-CellText $ws "D6" '152.65'
Set-CellText $ws "E6" '  +2.46%  '
Set-CellText $ws "D7" '3.519.05'
Set-CellText $ws "E7" '  +0.67%  '
Set-CellText $ws "E8" '  +0.11%  '
Set-CellText $ws "D9" '0.488'
Set-CellText $ws "E9" '  +1.75%  '
Set-CellText $ws "E10" '  +4.35%  '
Set-CellText $ws "D11" '7.60'
Set-CellText $ws "E11" '  +8.82%  '
Set-CellText $ws "D12" '0.434'
Set-CellText $ws "E12" '  +3.25%  '
Set-CellText $ws "D13" '32.94'
Set-CellText $ws "E13" '  +4.90%  '
Set-CellText $ws "D14" '0.0000217'
Set-CellText $ws "E14" '  +0.38%  '
Set-CellText $ws "D15" '4.132.86'
Set-CellText $ws "E15" '  +0.94%  '
Set-CellText $ws "D16" '3.533.90'
Set-CellText $ws "E16" '  +0.97%  '
Set-CellText $ws "D17" '68.148.46'
Set-CellText $ws "E17" '  +1.07%  '
Set-CellText $ws "E18" '  -0.48%  '
Set-CellText $ws "D19" '6.61'
Set-CellText $ws "E19" '  +3.58%  '
Set-CellText $ws "D20" '15.67'
Set-CellText $ws "E20" '  +3.98%  '
Set-CellText $ws "D21" '9.82'
Set-CellText $ws "E21" '  +8.33%  '
Set-CellText $ws "D22" '451.36'
Set-CellText $ws "E22" '  +1.41%  '
Set-CellText $ws "D23" '0.636'
Set-CellText $ws "E23" '  +2.80%  '
Set-CellText $ws "D24" '78.09'
Set-CellText $ws "E24" '  +1.13%  '
Set-CellText $ws "D25" '0.0000129'
Set-CellText $ws "E25" '  +2.46%  '
Set-CellText $ws "D26" '3.676.62'
Set-CellText $ws "E26" '  +0.95%  '
Set-CellText $ws "E27" '  -0.03%  '
Set-CellText $ws "D28" '9.03'
Set-CellText $ws "E28" '  +10.03%  '
Set-CellText $ws "D29" '10.16'
Set-CellText $ws "E29" '  +0.76%  '
Set-CellText $ws "D30" '1.67'
Set-CellText $ws "E30" '  +9.98%  '
Set-CellText $ws "D31" '2.52'
Set-CellText $ws "E31" '  +2.03%  '
Set-CellText $ws "E32" '  +4.07%  '
Set-CellText $ws "E33" '  -0.04%  '
Set-CellText $ws "D34" '25.84'
Set-CellText $ws "E34" '  +1.10%  '
Set-CellText $ws "D35" '6.24'
Set-CellText $ws "E35" '  +3.29%  '
Set-CellText $ws "D36" '1.89'
Set-CellText $ws "E36" '  +3.20%  '
Set-CellText $ws "D37" '3.526.86'
Set-CellText $ws "D38" '8.07'
Set-CellText $ws "E38" '  +0.59%  '
Set-CellText $ws "E39" '  +0.04%  '
Set-CellText $ws "D40" '2.34'
Set-CellText $ws "E40" '  +8.10%  '
Set-CellText $ws "D41" '1.00'
Set-CellText $ws "E41" '  +0.16%  '
Set-CellText $ws "D42" '0.0905'
Set-CellText $ws "E42" '  +3.23%  '
Set-CellText $ws "D43" '174.24'
Set-CellText $ws "E43" '  -2.21%  '
Set-CellText $ws "E44" '  +3.86%  '
Set-CellText $ws "D45" '30.91'
Set-CellText $ws "E45" '  +11.99%  '
Set-CellText $ws "D46" '0.884'
Set-CellText $ws "E46" '  +0.96%  '
Set-CellText $ws "D47" '46.76'
Set-CellText $ws "E47" '  +3.32%  '
Set-CellText $ws "D48" '1.32'
Set-CellText $ws "E48" '  +6.98%  '
Set-CellText $ws "D49" '2.56'
Set-CellText $ws "E49" '  +1.13%  '
Set-CellText $ws "D50" '7.71'
Set-CellText $ws "E50" '  +2.34%  '
Set-CellText $ws "D51" '0.257'
Set-CellText $ws "E51" '  +5.78%  '
